$d = $word.ActiveDocument
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>MADDE 1: TARAFLAR VE TANIMLAR</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:t>Bu</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> sözleşm</w:t>
      </w:r>
      <w:r>
        <w:t>e 1</w:t>
      </w:r>
      <w:r>
        <w:t>5</w:t>
      </w:r>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:r>
        <w:t>03/2021</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> tarihinde ticari merkezi </w:t>
      </w:r>
      <w:r>
        <w:t>Malatya</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_Hlk67072738"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>TEKNOLOJİLERİ</w:t>
      </w:r>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t xml:space="preserve">  ile</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> …………………………………… Adresinde</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>……………………………………(MÜŞTERİ olarak anılacaktır) arasında imzalanmıştır.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Anlaşma gereği MÜŞTERİ şirketimize </w:t>
      </w:r>
      <w:r>
        <w:t>300</w:t>
      </w:r>
      <w:r>
        <w:t>.000 (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>üçyüzbin</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">) Türk Lirası (TL) </w:t>
      </w:r>
      <w:r>
        <w:t>ödeyecektir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>MADDE 2: SÖZLEŞMENİN KONUSU</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:t>Bu</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> sözleşme </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>MÜŞTERİ’nin</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> iş süreçleri, kurum içi yönetim sistemleri, takip sistemleri vb. kısaca ticari faaliyetlerini sürdürürken ki ihtiyaç duyduğu tüm yazılım hizmetleri yürütebilmek amacıyla ihtiyaç duyduğu yazılımlarının </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>tarafından hazırlaması hizmetini kapsamaktadır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>MADDE 3: MÜŞTERİ’NİN YÜKÜMLÜLÜKLERİ</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>3.1.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> MÜŞTERİ, masa üstü ve web tabanlı yazılımların yükleneceği sunucu (sunucu / </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>hosting</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">) sistemini ve bu sistemin çalışabilmesi için gerekli donanım/yazılım gereksinimlerini kendisi tedarik etmekle yükümlüdür. Müşteri barındırma hizmetini ve alan adı tescilini </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>’den talep ederse o yılın ücret tablosu üzerinden hesaplanacak ve ek ücret sözleşme bedeline yansıtılacaktır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>3.2.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Yazılım süresi boyunca </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>TEKNOLOJİLERİ’in</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> talep ettiği </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>bilgi,evrak</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> ve yazılı dokümanları kendisi tedarik etmeli veya ettirmelidir. Müşterinin gerekli dokümanları zamanında vermemesinden kaynaklanan gecikmeden dolayı </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ sorumlu tutulamaz.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>3.3.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> MÜŞTERİ, </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ tarafından üretilen yazılımlarda kullanılan özel yazılım tekniklerinin telif hakkının </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ ‘e ait olduğunu, bu yazılımların hiçbir şekilde çoğaltılıp, dağıtılmayacağını kabul ve beyan eder.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>3.4.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> MÜŞTERİ, dile getirdiği tüm fikir, düşünce, ifade, yorum ve yazıların kendisine ait olduğunu, </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>TEKNOLOJİLERİ’</w:t>
      </w:r>
      <w:r>
        <w:t>n</w:t>
      </w:r>
      <w:r>
        <w:t>in</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> hiçbir şekilde sorumlu olmadığını kabul ve beyan eder.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">MADDE 4: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>KONTROL BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>’N</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>İN YÜKÜMLÜLÜKLERİ</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>4.1.</w:t>
      </w:r>
      <w:r>
        <w:t> Müşteri’nin ihtiyaç duyduğu yazılımları hazırlamakla yükümlüdür.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>4.2.</w:t>
      </w:r>
      <w:r>
        <w:t> Müşteri’ye sağlanmakta olan hizmetlerle ilgili oluşabilecek yazılım hatalarını ve sorunları 1 hafta içerisinde düzeltmekle yükümlüdür. Bu düzeltmeler için ücret talep edilmeyecektir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>4.3.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Müşteri için hazırlanmış tasarımdaki mevcut bölümlerdeki basit değişiklikler (komple tasarım ve kod/blok/div değişikliği gerektirmeyen düzeltmeler) için, gerekli çalışmayı tamamladıktan sonra </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>TEKNOLOJİLERİ’</w:t>
      </w:r>
      <w:r>
        <w:t>n</w:t>
      </w:r>
      <w:r>
        <w:t>in</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> belirleyeceği en kısa süre içinde gerekli düzenleme ve ilaveleri sisteme entegre edeceğini beyan eder. Bu düzeltmeler için ücret talep edilmeyecektir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>4.4.</w:t>
      </w:r>
      <w:r>
        <w:t> Müşteri tarafından talep edilen ek tasarım (tasarım ve kod değişikliği gerektiren düzeltmeler, sayfa eklemeler) için, gerekli çalışmayı tamamladıktan sonra belirleyeceği en kısa süre içinde gerekli düzenleme ve ilaveleri ek ücret mukabilinde sisteme entegre edeceğini beyan eder.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>4.4.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Müşteri tarafından iş süreçlerini </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>yazılımsal</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> olarak entegre etmekle yükümlü olup, istenilen işleyişin yerine getirmesiyle yükümlüdür.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>MADDE 5: GİZLİLİK VE GÜVENLİK</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">Hem </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>hem de MÜŞTERİ birbirlerinin gizlilik kurallarına saygı göstermelidir. Üçüncü şahıslarla paylaşılması zorunlu olmayan ve zaten diğer kişilerin ulaşımına açık olmayan tüm bilgiler gizli kabul edilip başka kişilerle paylaşılmamalıdır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>5.1.</w:t>
      </w:r>
      <w:r>
        <w:t> </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ, Dosya Transfer Yetkisi (FTP), Veri tabanı, Yönetim Kontrol Paneli şifreleri ile yazılım içeriğine ve veri tabanına kayıtlı bilgileri ve özel bilgileri 3. kişi ve kuruluşlarla paylaşmayacağını ve kullandırmayacağını beyan ve taahhüt eder.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>5.2.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> MÜŞTERİ, kendi çalışanlarının, herhangi bir kasıt, ihmal ya da kusurundan dolayı şifrelerin 3. kişi veya kuruluşların eline geçmesi halinde doğabilecek zararlardan </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>mesul değildir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>MADDE 6: ELEKTRONİK KAYITLARIN GEÇERLİLİĞİ</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>Taraflar arasında ve tarafların yetkililerince yapılan e-posta, anlık mesaj ve faks gibi elektronik yazışma kayıtları, kanunen geçerli delil sayılarak, usul hukuku bağlamında kesin ve bağlayıcı delil olarak kabul edilmiştir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>MADDE 7: TEBLİGATLAR</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:t>B</w:t>
      </w:r>
      <w:r>
        <w:t>u sözleşme ve uygulanması ile ilgili olarak yapılacak her türlü tebligat için, taraflar, işbu sözleşmede yazılı adresleri, yasal ikametgâh olarak belirlemişlerdir. Taraflar, bu adreslerde vaki değişiklikleri, diğer tarafa, noter kanalı veya iadeli taahhütlü mektup yolu ile bildirmedikleri takdirde, bu adreslerine yapılacak tebligatların geçerli, usulüne uygun ve kendilerine yapılmış sayılacağını kabul, beyan ve taahhüt ederler.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>MADDE 8: GARANTİ</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>Yazılımın, istenilen her türlü fonksiyonu yerine getirecek şekilde testleri tamamlandıktan sonra çalışır durumda yayımlandığı tarihten itibaren başlayacak ve toplam garanti süresi 1 yıl olacaktır.</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>tarafından titizlikle geliştirilecek ve test aşamasında olası hatalardan arındırılacak yazılımda, 1 yıllık garanti süresi içerisinde olabilecek yazılım veya tasarım hatalarının düzeltilmesinde veya yazılımın çalışmaması durumunda müdahale edilip çalışır duruma getirilmesi esnasında herhangi bir ücret talep edilmeyecektir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>MADDE 9: İŞİN TAMAMLANMASI &amp; İPTALİ</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">Sistemin hazırlanması Müşteri haricinde bir sebepten dolayı iş tamamlanmayacak ise </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>aldığı tüm ücretleri iade etmek zorundadır.</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">10.1 Müşteri için yapılmış olan tüm fonksiyonel ve tasarımsal ürünler/çalışmalar sözleşme sonrasında dijital olarak veya mail ile Müşteri’ye iletildiğinde / sunulduğunda </w:t>
      </w:r>
      <w:r>
        <w:t>KONTROL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> BİLİŞİM YAZILIM GELİŞTİRME TEKNOLOJİLERİ</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>işin %100’ünü tamamlamamış / teslim etmemiş olsa dahi peşin olarak aldığı sözleşme tutarında belirtilen %50’lik ön ödeme kısmını müşteri talep etmeyecektir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>MADDE 11: PROJE AMACI</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>Günümüz teknolojisine uygun gerek görsel gerekse de teknolojik olarak kurumu en iyi şekilde işleyişini sürdürecek ve yönetimsel marj sağlayacak müşterinin ihtiyaçları doğrultusu</w:t>
      </w:r>
      <w:r>
        <w:t>nd</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">a yazılımın müşteriye hazırlanıp </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>kapsamaktadır..</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>MADDE 1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>: YETKİLİ MAHKEME VE İCRA DAİRELERİ</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">İşbu sözleşmenin uygulanmasından doğabilecek her türlü uyuşmazlıkların çözümünde </w:t>
      </w:r>
      <w:r>
        <w:t>Malatya</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Mahkemeleri ve İcra Daireleri yetkilidir. 1</w:t>
      </w:r>
      <w:r>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (on </w:t>
      </w:r>
      <w:r>
        <w:t>bir)</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> maddeden oluşan bu sözleşm</w:t>
      </w:r>
      <w:r>
        <w:t>e 13</w:t>
      </w:r>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:r>
        <w:t>03/2021</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> tarihinde 2 (İki) nüsha olarak düzenlenmiş olup, okunmuş ve altına imza alınmıştır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t> </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
